$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object "object[,]" 24,5
$arrBF[0,0] = 19.54007855928528
$arrBF[0,1] = 6.403381025117031
$arrBF[0,2] = 4.736560110516001
$arrBF[0,3] = 10.65486524604357
$arrBF[0,4] = 50.32538137073519
$arrBF[1,0] = 19.42173598948629
$arrBF[1,1] = 6.274015419378872
$arrBF[1,2] = 4.735288515693775
$arrBF[1,3] = 10.67026975130645
$arrBF[1,4] = 50.34835583333942
$arrBF[2,0] = 19.35328567392238
$arrBF[2,1] = 6.195227297311263
$arrBF[2,2] = 4.735032054012245
$arrBF[2,3] = 10.68076458328736
$arrBF[2,4] = 50.37157890864855
$arrBF[3,0] = 19.32647562647852
$arrBF[3,1] = 6.163331557302097
$arrBF[3,2] = 4.735060107510242
$arrBF[3,3] = 10.68530240393107
$arrBF[3,4] = 50.38333415167596
$arrBF[4,0] = 19.32208999918373
$arrBF[4,1] = 6.158049492849307
$arrBF[4,2] = 4.735072792810984
$arrBF[4,3] = 10.68607168876062
$arrBF[4,4] = 50.38542447128061
$arrBF[5,0] = 19.35291968374561
$arrBF[5,1] = 6.194796220061566
$arrBF[5,2] = 4.735031894662402
$arrBF[5,3] = 10.68082472413879
$arrBF[5,4] = 50.37172816692282
$arrBF[6,0] = 19.49841510008703
$arrBF[6,1] = 6.358671809577586
$arrBF[6,2] = 4.736013288292538
$arrBF[6,3] = 10.65996190368337
$arrBF[6,4] = 50.3314105432867
$arrBF[7,0] = 19.81598962260647
$arrBF[7,1] = 6.682997188148571
$arrBF[7,2] = 4.742065288817394
$arrBF[7,3] = 10.62725225823497
$arrBF[7,4] = 50.32469321505887
$arrBF[8,0] = 20.06731244532062
$arrBF[8,1] = 6.920312362171356
$arrBF[8,2] = 4.748982927110503
$arrBF[8,3] = 10.60819232959844
$arrBF[8,4] = 50.36383344612489
$arrBF[9,0] = 20.18517289337994
$arrBF[9,1] = 7.027509599923006
$arrBF[9,2] = 4.752656299845666
$arrBF[9,3] = 10.60059485629832
$arrBF[9,4] = 50.39118621767989
$arrBF[10,0] = 20.23027788076852
$arrBF[10,1] = 7.067950097207169
$arrBF[10,2] = 4.75412208700344
$arrBF[10,3] = 10.59787165453562
$arrBF[10,4] = 50.40291318713766
$arrBF[11,0] = 20.22054319497171
$arrBF[11,1] = 7.059247976453714
$arrBF[11,2] = 4.753803094011109
$arrBF[11,3] = 10.59845131265801
$arrBF[11,4] = 50.40032675650187
$arrBF[12,0] = 20.18887441082524
$arrBF[12,1] = 7.030839967208931
$arrBF[12,2] = 4.75277539828642
$arrBF[12,3] = 10.60036773742681
$arrBF[12,4] = 50.39212358990678
$arrBF[13,0] = 20.16953704035953
$arrBF[13,1] = 7.013418091400752
$arrBF[13,2] = 4.752155612239423
$arrBF[13,3] = 10.60156161683077
$arrBF[13,4] = 50.38727707546752
$arrBF[14,0] = 20.05967816775937
$arrBF[14,1] = 6.913287715875863
$arrBF[14,2] = 4.748753377754144
$arrBF[14,3] = 10.60871039126755
$arrBF[14,4] = 50.36223769211857
$arrBF[15,0] = 19.99316400978024
$arrBF[15,1] = 6.851634830275158
$arrBF[15,2] = 4.746800402513376
$arrBF[15,3] = 10.61337039939147
$arrBF[15,4] = 50.34931977350779
$arrBF[16,0] = 19.95524141836948
$arrBF[16,1] = 6.816104780693736
$arrBF[16,2] = 4.745726718975305
$arrBF[16,3] = 10.61615173757435
$arrBF[16,4] = 50.34278869305945
$arrBF[17,0] = 19.94245997449503
$arrBF[17,1] = 6.804064446268783
$arrBF[17,2] = 4.745371738171531
$arrBF[17,3] = 10.61711081696006
$arrBF[17,4] = 50.34073188226687
$arrBF[18,0] = 20.00021017119134
$arrBF[18,1] = 6.858205344400788
$arrBF[18,2] = 4.74700317175788
$arrBF[18,3] = 10.6128638814606
$arrBF[18,4] = 50.35060189353506
$arrBF[19,0] = 20.19816371798312
$arrBF[19,1] = 7.039188570256455
$arrBF[19,2] = 4.753075236193998
$arrBF[19,3] = 10.59980066704343
$arrBF[19,4] = 50.39449593822697
$arrBF[20,0] = 20.33028306839995
$arrBF[20,1] = 7.156562603624883
$arrBF[20,2] = 4.757478978574115
$arrBF[20,3] = 10.59215930631887
$arrBF[20,4] = 50.43116074818019
$arrBF[21,0] = 20.25952877586883
$arrBF[21,1] = 7.094014835495543
$arrBF[21,2] = 4.755089111370555
$arrBF[21,3] = 10.59615580867614
$arrBF[21,4] = 50.41086363071562
$arrBF[22,0] = 19.99702361068582
$arrBF[22,1] = 6.855235076325024
$arrBF[22,2] = 4.746911346765055
$arrBF[22,3] = 10.61309255973048
$arrBF[22,4] = 50.35001945691521
$arrBF[23,0] = 19.72679021536637
$arrBF[23,1] = 6.595233653642055
$arrBF[23,2] = 4.73999069228355
$arrBF[23,3] = 10.63522587763115
$arrBF[23,4] = 50.31876473596024
$ws.Range("B2:F25").Value = $arrBF

$arrIL = New-Object "object[,]" 24,4
$arrIL[0,0] = 36.22301452445411
$arrIL[0,1] = 9.96456123039037
$arrIL[0,2] = 17.35877020247769
$arrIL[0,3] = 12.05623818284223
$arrIL[1,0] = 36.27675806777527
$arrIL[1,1] = 9.98121583488812
$arrIL[1,2] = 17.28004288096013
$arrIL[1,3] = 12.06628163864897
$arrIL[2,0] = 36.31531285776769
$arrIL[2,1] = 9.992088785473468
$arrIL[2,2] = 17.23540360998659
$arrIL[2,3] = 12.07415306756729
$arrIL[3,0] = 36.33242004420154
$arrIL[3,1] = 9.996682704421925
$arrIL[3,2] = 17.21815753813523
$arrIL[3,3] = 12.07778988433812
$arrIL[4,0] = 36.33534493986817
$arrIL[4,1] = 9.997455384599499
$arrIL[4,2] = 17.21535130294288
$arrIL[4,3] = 12.07841970682214
$arrIL[5,0] = 36.31553792168661
$arrIL[5,1] = 9.992150079713007
$arrIL[5,2] = 17.23516717981155
$arrIL[5,3] = 12.07420037678045
$arrIL[6,0] = 36.24039140489439
$arrIL[6,1] = 9.970169715456839
$arrIL[6,2] = 17.33086623704547
$arrIL[6,3] = 12.05934766689048
$arrIL[7,0] = 36.13716800479834
$arrIL[7,1] = 9.932180642495442
$arrIL[7,2] = 17.54719325208671
$arrIL[7,3] = 12.0437242430882
$arrIL[8,0] = 36.08830462229676
$arrIL[8,1] = 9.90736175655587
$arrIL[8,2] = 17.72259546934974
$arrIL[8,3] = 12.04044265439569
$arrIL[9,0] = 36.07194284718393
$arrIL[9,1] = 9.896736805067905
$arrIL[9,2] = 17.80572971091572
$arrIL[9,3] = 12.04072051087905
$arrIL[10,0] = 36.06659100947497
$arrIL[10,1] = 9.892808656703867
$arrIL[10,2] = 17.83766983722227
$arrIL[10,3] = 12.04107938619877
$arrIL[11,0] = 36.06770608170396
$arrIL[11,1] = 9.893650421645125
$arrIL[11,2] = 17.83077089980884
$arrIL[11,3] = 12.04099083089675
$arrIL[12,0] = 36.07148563147753
$arrIL[12,1] = 9.896411726326949
$arrIL[12,2] = 17.8083483726083
$arrIL[12,3] = 12.0407449579687
$arrIL[13,0] = 36.07391063843625
$arrIL[13,1] = 9.898115503644558
$arrIL[13,2] = 17.79467304371494
$arrIL[13,3] = 12.04062735707663
$arrIL[14,0] = 36.08949199240292
$arrIL[14,1] = 9.908069477043371
$arrIL[14,2] = 17.71722788146855
$arrIL[14,3] = 12.0404600461472
$arrIL[15,0] = 36.10055356297321
$arrIL[15,1] = 9.914346043261665
$arrIL[15,2] = 17.67055879427295
$arrIL[15,3] = 12.04081026324037
$arrIL[16,0] = 36.10746805119965
$arrIL[16,1] = 9.91801880013902
$arrIL[16,2] = 17.64403197321814
$arrIL[16,3] = 12.04117842632747
$arrIL[17,0] = 36.10990399602328
$arrIL[17,1] = 9.919273102879457
$arrIL[17,2] = 17.63510536935961
$arrIL[17,3] = 12.04133174334085
$arrIL[18,0] = 36.0993188889091
$arrIL[18,1] = 9.913671411600694
$arrIL[18,2] = 17.67549424776279
$arrIL[18,3] = 12.04075573302323
$arrIL[19,0] = 36.07035257770094
$arrIL[19,1] = 9.895598081044865
$arrIL[19,2] = 17.81492213240165
$arrIL[19,3] = 12.0408103007917
$arrIL[20,0] = 36.05634098336435
$arrIL[20,1] = 9.884341366441765
$arrIL[20,2] = 17.90871106935071
$arrIL[20,3] = 12.04232395488627
$arrIL[21,0] = 36.06336901854228
$arrIL[21,1] = 9.89029860490766
$arrIL[21,2] = 17.85841769647781
$arrIL[21,3] = 12.04138119667434
$arrIL[22,0] = 36.09987535606176
$arrIL[22,1] = 9.913976212330759
$arrIL[22,2] = 17.67326198107827
$arrIL[22,3] = 12.04077986643751
$arrIL[23,0] = 36.16035954471083
$arrIL[23,1] = 9.941912884130458
$arrIL[23,2] = 17.48570213262312
$arrIL[23,3] = 12.04650844682682
$ws.Range("I2:L25").Value = $arrIL

$arrN = New-Object "object[,]" 24,1
$arrN[0,0] = 24.47797805936221
$arrN[1,0] = 24.53379105072677
$arrN[2,0] = 24.56999178885812
$arrN[3,0] = 24.5852302739086
$arrN[4,0] = 24.58779001148946
$arrN[5,0] = 24.57019533006599
$arrN[6,0] = 24.49682197577677
$arrN[7,0] = 24.36822724363858
$arrN[8,0] = 24.2830218786368
$arrN[9,0] = 24.24626333430226
$arrN[10,0] = 24.23263099853324
$arrN[11,0] = 24.23555419521969
$arrN[12,0] = 24.24513603833628
$arrN[13,0] = 24.25104259700333
$arrN[14,0] = 24.28546437163708
$arrN[15,0] = 24.30709334239473
$arrN[16,0] = 24.31972220385936
$arrN[17,0] = 24.32403050224868
$arrN[18,0] = 24.30477140157489
$arrN[19,0] = 24.24231382489119
$arrN[20,0] = 24.20316871881346
$arrN[21,0] = 24.22390814899569
$arrN[22,0] = 24.30582054637197
$arrN[23,0] = 24.40138365050888
$ws.Range("N2:N25").Value = $arrN

